# Auto-generated edit script: updates cryptos list values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.009.70"
$ws.Range("E2").Value = "  +6.28%  "
$ws.Range("D3").Value = "'2.259.04"
$ws.Range("E3").Value = "  +4.28%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "'233.64"
$ws.Range("E5").Value = "  +2.09%  "
$ws.Range("D6").Value = "'0.641"
$ws.Range("E6").Value = "  +1.50%  "
$ws.Range("D7").Value = "'62.83"
$ws.Range("E7").Value = "  -1.48%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "'0.409"
$ws.Range("E9").Value = "  +3.43%  "
$ws.Range("D10").Value = "'59.69"
$ws.Range("E10").Value = "  +2.84%  "
$ws.Range("D11").Value = "'0.0894"
$ws.Range("E11").Value = "  +4.95%  "
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("D13").Value = "'2.602.56"
$ws.Range("E13").Value = "  +4.60%  "
$ws.Range("D14").Value = "'15.95"
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").Value = "'22.76"
$ws.Range("E15").Value = "  +3.15%  "
$ws.Range("D16").Value = "'0.819"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").Value = "'5.68"
$ws.Range("E17").Value = "  +2.92%  "
$ws.Range("D18").Value = "'2.262.99"
$ws.Range("E18").Value = "  +4.59%  "
$ws.Range("D19").Value = "'41.999.04"
$ws.Range("E19").Value = "  +6.23%  "
$ws.Range("D20").Value = "'74.40"
$ws.Range("E20").Value = "  +3.06%  "
$sub3 = [char]0x2083
$ws.Range("D21").Value = "'0.0${sub3}0919"
$ws.Range("E21").Value = "  +8.33%  "
$ws.Range("D22").Value = "'6.12"
$ws.Range("E22").Value = "  -1.53%  "
$ws.Range("D23").Value = "'251.24"
$ws.Range("E23").Value = "  +9.17%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'2.42"
$ws.Range("E25").Value = "  +3.59%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'2.41"
$ws.Range("E26").Value = "  +2.40%  "
$ws.Range("D27").Value = "'0.149"
$ws.Range("E27").Value = "  +7.17%  "
$ws.Range("D28").Value = "'9.80"
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("D29").Value = "'170.86"
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'1.48"
$ws.Range("E30").Value = "  +4.24%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'20.42"
$ws.Range("E31").Value = "  +2.58%  "
$ws.Range("E32").Value = "  +6.60%  "
$ws.Range("D33").Value = "'0.124"
$ws.Range("E33").Value = "  +1.25%  "
$ws.Range("D34").Value = "'5.05"
$ws.Range("E34").Value = "  +7.22%  "
$ws.Range("D35").Value = "'4.77"
$ws.Range("E35").Value = "  +3.43%  "
$ws.Range("D36").Value = "'0.0634"
$ws.Range("E36").Value = "  +2.12%  "
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").Value = "'6.74"
$ws.Range("E37").Value = "  -4.16%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'3.78"
$ws.Range("E38").Value = "  +2.60%  "
$ws.Range("D39").Value = "'2.44"
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("B40").Value = "BinanceUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D40").Value = "'1.01"
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("D41").Value = "'5.03"
$ws.Range("E41").Value = "  +14.89%  "
$ws.Range("B42").Value = "TerraClassic"
$ws.Range("C42").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D42").Value = "'0.000241"
$ws.Range("E42").Value = "  +31.33%  "
$ws.Range("D43").Value = "'0.0241"
$ws.Range("E43").Value = "  +5.28%  "
$ws.Range("D44").Value = "'8.62"
$ws.Range("E44").Value = "  +10.71%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'100.97"
$ws.Range("E45").Value = "  -1.74%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'17.50"
$ws.Range("E46").Value = "  -2.57%  "
$ws.Range("D47").Value = "'1.23"
$ws.Range("E47").Value = "  +2.00%  "
$ws.Range("D48").Value = "'0.0976"
$ws.Range("E48").Value = "  +5.49%  "
$ws.Range("D49").Value = "'1.499.42"
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("D50").Value = "'1.12"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("E51").Value = "  -0.38%  "
